$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 120
$ws.Range("D120").Value = 44642
$ws.Range("J120").Value = 220
$ws.Range("K120").Value = 17000
$ws.Range("L120").Value = 18000
$ws.Range("M120").Value = 17455
$ws.Range("P120").Value = 1746

# Row 121
$ws.Range("D121").Value = 44357
$ws.Range("J121").Value = 400
$ws.Range("K121").Value = 12000
$ws.Range("L121").Value = 13000
$ws.Range("M121").Value = 12500
$ws.Range("P121").Value = 1250

# Row 122
$ws.Range("D122").Value = 44537
$ws.Range("J122").Value = 310
$ws.Range("K122").Value = 17000
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = 17516
$ws.Range("P122").Value = 1752

# Row 123
$ws.Range("D123").Value = 44490
$ws.Range("K123").Value = 16000
$ws.Range("L123").Value = 16500
$ws.Range("M123").Value = 16250
$ws.Range("P123").Value = 1625

# Row 124
$ws.Range("D124").Value = 44397
$ws.Range("K124").Value = 13000
$ws.Range("L124").Value = 14000
$ws.Range("M124").Value = 13500
$ws.Range("P124").Value = 1350

# Row 125
$ws.Range("D125").Value = 44462

# Row 126
$ws.Range("D126").Value = 44446
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 15500
$ws.Range("M126").Value = 15250
$ws.Range("P126").Value = 1525

# Row 127
$ws.Range("D127").Value = 44329
$ws.Range("J127").Value = 400
$ws.Range("K127").Value = 11000
$ws.Range("L127").Value = 12000
$ws.Range("M127").Value = 11500
$ws.Range("P127").Value = 1150

# Row 128
$ws.Range("D128").Value = 44637
$ws.Range("J128").Value = 150
$ws.Range("K128").Value = 17000
$ws.Range("L128").Value = 18000
$ws.Range("M128").Value = 17467
$ws.Range("P128").Value = 1747

# Row 129
$ws.Range("D129").Value = 44208
$ws.Range("J129").Value = 300
$ws.Range("K129").Value = 11000
$ws.Range("L129").Value = 12000
$ws.Range("M129").Value = 11333
$ws.Range("P129").Value = 1133

# Row 130
$ws.Range("D130").Value = 44355
$ws.Range("J130").Value = 400
$ws.Range("K130").Value = 12000
$ws.Range("L130").Value = 13000
$ws.Range("M130").Value = 12500
$ws.Range("P130").Value = 1250

# Row 131
$ws.Range("I131").Value = "Primera"
$ws.Range("J131").Value = 310
$ws.Range("K131").Value = 21000
$ws.Range("L131").Value = 22000
$ws.Range("M131").Value = 21484
$ws.Range("P131").Value = 2148

# Row 132
$ws.Range("D132").Value = 44530
$ws.Range("I132").Value = "Segunda"
$ws.Range("J132").Value = 290
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 19000
$ws.Range("M132").Value = 18483
$ws.Range("P132").Value = 1848

# Row 133
$ws.Range("D133").Value = 44483
$ws.Range("J133").Value = 450
$ws.Range("K133").Value = 15000
$ws.Range("L133").Value = 16000
$ws.Range("M133").Value = 15444
$ws.Range("P133").Value = 1544

# Row 134
$ws.Range("D134").Value = 44294

# Row 135
$ws.Range("D135").Value = 44264
$ws.Range("K135").Value = 12000
$ws.Range("L135").Value = 13000
$ws.Range("M135").Value = 12500
$ws.Range("P135").Value = 1250

# Row 136
$ws.Range("D136").Value = 44232
$ws.Range("K136").Value = 11000
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = 11500
$ws.Range("P136").Value = 1150

# Row 137
$ws.Range("D137").Value = 44330
$ws.Range("J137").Value = 400
$ws.Range("K137").Value = 12000
$ws.Range("L137").Value = 13000
$ws.Range("M137").Value = 12500
$ws.Range("P137").Value = 1250

# Row 138
$ws.Range("D138").Value = 44504
$ws.Range("I138").Value = "Primera"
$ws.Range("J138").Value = 350
$ws.Range("K138").Value = 14000
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = 14571
$ws.Range("N138").Value = "`$/caja 10 kilos"
$ws.Range("O138").Value = "China"
$ws.Range("P138").Value = 1457
$ws.Range("Q138").Value = 10

# Row 139
$ws.Range("I139").Value = "2a (cosecha)"
$ws.Range("J139").Value = 310
$ws.Range("K139").Value = 6500
$ws.Range("L139").Value = 7000
$ws.Range("M139").Value = 6742
$ws.Range("N139").Value = "`$/trenza 50 unidades"
$ws.Range("O139").Value = "Provincia de Talagante"
$ws.Range("P139").Value = 1348
$ws.Range("Q139").Value = 5

# Row 140
$ws.Range("D140").Value = 44572
$ws.Range("J140").Value = 270
$ws.Range("K140").Value = 14000
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = 14444
$ws.Range("P140").Value = 1444

# Row 141
$ws.Range("D141").Value = 44257

# Row 142
$ws.Range("D142").Value = 44370
$ws.Range("K142").Value = 12000
$ws.Range("M142").Value = 12500
$ws.Range("P142").Value = 1250

# Row 143
$ws.Range("D143").Value = 44385
$ws.Range("K143").Value = 12500
$ws.Range("M143").Value = 12750
$ws.Range("P143").Value = 1275

# Row 144
$ws.Range("D144").Value = 44236
$ws.Range("K144").Value = 12000
$ws.Range("L144").Value = 13000
$ws.Range("M144").Value = 12500
$ws.Range("P144").Value = 1250

# Row 145
$ws.Range("D145").Value = 44229
$ws.Range("K145").Value = 11000
$ws.Range("L145").Value = 12000
$ws.Range("M145").Value = 11500
$ws.Range("P145").Value = 1150

# Row 146
$ws.Range("D146").Value = 44299
$ws.Range("K146").Value = 14000
$ws.Range("L146").Value = 15000
$ws.Range("M146").Value = 14500
$ws.Range("P146").Value = 1450

# Row 147
$ws.Range("D147").Value = 44610
$ws.Range("K147").Value = 17000
$ws.Range("L147").Value = 18000
$ws.Range("M147").Value = 17500
$ws.Range("P147").Value = 1750

# Row 148
$ws.Range("D148").Value = 44399
$ws.Range("L148").Value = 13500
$ws.Range("M148").Value = 13250
$ws.Range("P148").Value = 1325

# Row 149
$ws.Range("D149").Value = 44390
$ws.Range("K149").Value = 13000
$ws.Range("L149").Value = 14000
$ws.Range("M149").Value = 13500
$ws.Range("P149").Value = 1350

# Row 150
$ws.Range("D150").Value = 44285
$ws.Range("J150").Value = 400
$ws.Range("K150").Value = 12000
$ws.Range("L150").Value = 13000
$ws.Range("M150").Value = 12500
$ws.Range("P150").Value = 1250

# Row 151
$ws.Range("D151").Value = 44498
$ws.Range("J151").Value = 350
$ws.Range("M151").Value = 17571
$ws.Range("P151").Value = 1757

# Row 152
$ws.Range("D152").Value = 44595
$ws.Range("J152").Value = 400
$ws.Range("M152").Value = 17500
$ws.Range("P152").Value = 1750

# Row 153
$ws.Range("D153").Value = 44628
$ws.Range("J153").Value = 220
$ws.Range("M153").Value = 17455
$ws.Range("P153").Value = 1746

# Row 154
$ws.Range("A154").Value = 11
$ws.Range("B154").Value = "Vega Monumental Concepción"
$ws.Range("C154").Value = "Bíobío"
$ws.Range("D154").Value = 44552
$ws.Range("E154").Value = 8
$ws.Range("F154").Value = 100112003
$ws.Range("G154").Value = "Ajo"
$ws.Range("H154").Value = "Chino"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 400
$ws.Range("K154").Value = 17000
$ws.Range("L154").Value = 18000
$ws.Range("M154").Value = 17500
$ws.Range("N154").Value = "`$/caja 10 kilos"
$ws.Range("O154").Value = "China"
$ws.Range("P154").Value = 1750
$ws.Range("Q154").Value = 10
$ws.Range("R154").Value = "Hortaliza"

# Ensure D154 uses the same date number format as the rest of column D
$ws.Range("D154").NumberFormat = "YYYY-MM-DD HH:MM:SS"
